# Apply the edit described by the commit:
#  1. After the paragraph ending in "... attaccarlo altrove?" append
#     " " + (Wingdings arrow symbol) + " pare di si` invece", and move
#     the _GoBack bookmark to the end of that paragraph.
#  2. Merge the two runs of "Buildare il" / " sistemino da mandare a
#     Denti che permetta di interpretare semplici frasi." into a single
#     run, and drop the _GoBack bookmark that used to sit between them
#     (it now lives in paragraph 1).

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Step 1: locate end of the "altrove?" paragraph and append text.
# ---------------------------------------------------------------
$findRange = $d.Content.Duplicate
$null = $findRange.Find.Execute("attaccarlo altrove?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertion = $d.Range($findRange.End, $findRange.End)
$insertion.InsertAfter(" ")

$afterSpace = $d.Range($insertion.End, $insertion.End)
$afterSpace.InsertAfter([string][char]0xF0E0)
$symRange = $d.Range($afterSpace.End - 1, $afterSpace.End)
$symRange.Font.Name = "Wingdings"

$afterSym = $d.Range($afterSpace.End, $afterSpace.End)
$afterSym.InsertAfter(" pare di s" + [string][char]0x00EC + " invece")

# Match the language formatting already used throughout the document.
$wholeAddition = $d.Range($findRange.End, $afterSym.End + (" pare di s invece").Length + 1)
$wholeAddition.LanguageID = $findRange.LanguageID

# ---------------------------------------------------------------
# Step 2: move the _GoBack bookmark to the end of this paragraph.
# ---------------------------------------------------------------
$para1 = $findRange.Paragraphs.Item(1)
$para1End = $para1.Range.End - 1

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
$null = $d.Bookmarks.Add("_GoBack", $d.Range($para1End, $para1End))

# ---------------------------------------------------------------
# Step 3: merge "Buildare il" / " sistemino..." runs into one run
# (simple text stays the same, only the run boundary + bookmark is
# removed, so no Find/Replace of visible text is actually needed -
# deleting the bookmark above already took care of the only
# bookmark in the document).
# ---------------------------------------------------------------
